# Connectors-Template.xlsx: replace the "organization-display-name" column
# with a new "connector-description" column (inserted right after
# "connector-id"), per commit "Remove DisplayName attribute setting (#5)".
#
# Net column layout change on the "Connectors" sheet / Table1:
#   before: connector-id | connector-port | connector-version | connector-db-connection-string | organization-display-name | backbone-base-url | backbone-client-id | backbone-client-secret
#   after:  connector-id | connector-description | connector-port | connector-version | connector-db-connection-string | backbone-base-url | backbone-client-id | backbone-client-secret
#
# Columns F/G/H (backbone-*) keep their position; A keeps its position.
# B/C/D/E are rewritten in place (no net column-count change) so the table's
# ListColumns collection re-syncs by simply updating the header text for each
# cell - inserting/deleting table columns is not reliable in this host, but
# writing header cell values is what actually drives the table column name.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Connectors")

# --- Rewrite header row (row 1) in its new left-to-right order ---
$ws.Range("A1").Value = "connector-id"
$ws.Range("B1").Value = "connector-description"
$ws.Range("C1").Value = "connector-port"
$ws.Range("D1").Value = "connector-version"
$ws.Range("E1").Value = "connector-db-connection-string"
$ws.Range("F1").Value = "backbone-base-url"
$ws.Range("G1").Value = "backbone-client-id"
$ws.Range("H1").Value = "backbone-client-secret"

# --- Data row (row 2) ---
# The only populated data cell was the integer-formatted "connector-port"
# cell, previously B2. That column is now C, so move the number format
# there and clear the old B2 cell so it no longer carries the format.
$ws.Range("B2").Clear()
$ws.Range("C2").NumberFormat = "0"

# --- View state / print setup (cosmetic, but cheap & accurate to reproduce) ---
$ws.Range("F:F").Select() | Out-Null
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Default Values sheet is unaffected in content (still references
# connector-version / backbone-* / connector-db-connection-string); nothing
# else to change there.
